# rider-template.xlsx: add gsm/username/religion/gender columns to the
# roster sheet (renamed Sheet1 -> Sheet2) and append a new "NOTE" sheet
# that explains the numeric codes used for gender/religion.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Roster sheet: rename Sheet1 -> Sheet2, add the new columns
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# re-type the columns that stay put first (so the shared-string table
# keeps its original ordering), then the shifted / brand-new ones
$ws1.Range("A1").Value = "firstname"
$ws1.Range("B1").Value = "surname"
$ws1.Range("C1").Value = "othername"
$ws1.Range("F1").Value = "email"
$ws1.Range("G1").Value = "address"
$ws1.Range("H1").Value = "religion"
$ws1.Range("I1").Value = "gender"
$ws1.Range("D1").Value = "gsm"
$ws1.Range("E1").Value = "username"

# header emphasis (red text) on a few columns
$ws1.Range("A1,B1,D1,G1").Font.Color = 255

$ws1.Columns.Item(1).ColumnWidth = 9.57
$ws1.Columns.Item(2).ColumnWidth = 8.71
$ws1.Columns.Item(3).ColumnWidth = 10.86
$ws1.Columns.Item(4).ColumnWidth = 11
$ws1.Columns.Item(5).ColumnWidth = 9.86

# freeze the header row
$ws1.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# data validation rules
$dvGender = $ws1.Range("I1:I1048576")
$dvGender.Validation.Add(1, 1, 1, 1, 3)
$dvGender.Validation.ErrorMessage = "use 1 for Male, 2 for Female and 3 for other"

$dvReligion = $ws1.Range("H1:H1048576")
$dvReligion.Validation.Add(1, 1, 1, 1, 3)
$dvReligion.Validation.ErrorMessage = "use 1 for Muslim, 2 for Christain and 3 for Other"

$dvPhone = $ws1.Range("D1:D1048576")
$dvPhone.Validation.Add(6, 1, 3, 10)
$dvPhone.Validation.ErrorMessage = "10 digits, omit the first 0"
$dvPhone.Validation.InputTitle = "phone number"
$dvPhone.Validation.IgnoreBlank = $false

$dvSurname = $ws1.Range("B1:B1048576")
$dvSurname.Validation.Add(6, 1, 6, 30)
$dvSurname.Validation.ErrorMessage = "maximum of 30 character"
$dvSurname.Validation.IgnoreBlank = $false

$dvFirstname = $ws1.Range("A1:A1048576")
$dvFirstname.Validation.Add(6, 1, 6, 31)
$dvFirstname.Validation.ErrorMessage = "maximum of 30 character"

$dvAddress = $ws1.Range("G1")
$dvAddress.Validation.Add(6, 1, 6, 256)
$dvAddress.Validation.ErrorMessage = "address is compulsary and not more than 255 characters"

# select H8 as the final cursor position (matches the saved view)
$ws1.Range("H8").Select()
$ws1.Name = "Sheet2"

# ---------------------------------------------------------------------
# 2. New NOTE sheet explaining the gender/religion codes
# ---------------------------------------------------------------------
$note = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$note.Name = "NOTE"

$note.Range("A1").Value = "GENDER "
$note.Range("A2").Value = "the number coresponding number represent each gender"
$note.Range("A3").Value = "MALE"
$note.Range("B3").Value = 1
$note.Range("A4").Value = "FEMALE"
$note.Range("B4").Value = 2
$note.Range("A5").Value = "OTHER"
$note.Range("B5").Value = 3

$note.Range("D1").Value = "RELIGION"
$note.Range("D2").Value = "the corresponding number represent each religion"
$note.Range("D3").Value = "MUSLIM"
$note.Range("E3").Value = 1
$note.Range("D4").Value = "CHRISTIAN"
$note.Range("E4").Value = 2
$note.Range("D5").Value = "OTHER"
$note.Range("E5").Value = 3

# -- title row 1 : big red, centered, wrapped, merged --
$note.Range("A1:B1").Font.Size = 18
$note.Range("A1:B1").Font.Color = 255
$note.Range("A1:B1").HorizontalAlignment = -4108
$note.Range("A1:B1").WrapText = $true
$note.Range("A1:B1").Merge()

$note.Range("D1:E1").Font.Size = 16
$note.Range("D1:E1").Font.Color = 255
$note.Range("D1:E1").HorizontalAlignment = -4108
$note.Range("D1:E1").Merge()

# -- explanatory row 2 : centered, wrapped, bottom border, merged --
$note.Range("A2:B2").Font.Size = 18
$note.Range("A2:B2").HorizontalAlignment = -4108
$note.Range("A2:B2").WrapText = $true
$note.Range("A2:B2").Borders.Item(9).LineStyle = 1
$note.Range("A2:B2").Merge()

$note.Range("D2:E2").Font.Size = 16
$note.Range("D2:E2").HorizontalAlignment = -4108
$note.Range("D2:E2").WrapText = $true
$note.Range("D2:E2").Borders.Item(9).LineStyle = 1
$note.Range("D2:E2").Merge()

# -- code blocks rows 3-5 : white-on-colour, boxed --
$green = 5287936
$darkred = 192

$note.Range("A3:B5").Font.Size = 48
$note.Range("A3:B5").Font.ThemeColor = 2
$note.Range("A3:B5").Interior.Color = $green
$note.Range("A3:B5").Borders.LineStyle = 1
$note.Range("B3:B5").HorizontalAlignment = -4108

$note.Range("D3:E5").Font.Size = 48
$note.Range("D3:E5").Font.ThemeColor = 2
$note.Range("D3:E5").Interior.Color = $darkred
$note.Range("D3:E5").Borders.LineStyle = 1
$note.Range("E3:E5").HorizontalAlignment = -4108
$note.Range("E3:E5").VerticalAlignment = -4108

# -- row heights --
$note.Rows.Item(1).RowHeight = 23.25
$note.Rows.Item(2).RowHeight = 23.25
$note.Rows.Item(3).RowHeight = 61.5
$note.Rows.Item(4).RowHeight = 61.5
$note.Rows.Item(5).RowHeight = 61.5
$note.Rows.Item(8).RowHeight = 23.25

# -- trailing merged, empty, bordered block on row 8 --
$note.Range("D8:E8").Font.Size = 18
$note.Range("D8:E8").Font.Color = 255
$note.Range("D8:E8").Borders.Item(9).LineStyle = 1
$note.Range("D8:E8").Merge()

# -- column widths --
$note.Columns.Item(1).ColumnWidth = 33.43
$note.Columns.Item(2).ColumnWidth = 17.29
$note.Columns.Item(4).ColumnWidth = 43.57
$note.Columns.Item(5).ColumnWidth = 11.14

$note.Range("B13").Select()
$ws1.Select()

Write-Output "applied rider-template edits"
